# The "(minimization) target" and "(optimization) target" labels (and a
# stray double-space variant "(minimization)  target") are being unified
# into a single "target" label throughout the sheet (comment parsing -
# not yet serialized). Update every cell in column B that currently
# holds one of those variants so it just reads "target".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(24, 27, 36, 39, 46, 49, 58, 61, 65, 68, 78, 81)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "target"
}
